$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.230.14'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.085.98'
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.82'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.34'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.077.85'
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("E10").Value = '  +5.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.60'
$ws.Range("E11").Value = '  -2.15%  '
$ws.Range("E12").Value = '  -2.78%  '
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("E14").Value = '  +5.72%  '
$ws.Range("E15").Value = '  -1.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.595.95'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.111.64'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.09'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.086.22'
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '460.36'
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.20'
$ws.Range("E21").Value = '  +0.85%  '
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.42'
$ws.Range("E23").Value = '  -1.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.94'
$ws.Range("E24").Value = '  -2.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.07'
$ws.Range("E25").Value = '  -1.28%  '
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.95'
$ws.Range("E28").Value = '  +8.15%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.20'
$ws.Range("E30").Value = '  -2.07%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.66'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.81'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.110'
$ws.Range("E33").Value = '  -1.23%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.55'
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("E35").Value = '  -1.52%  '
$ws.Range("E36").Value = '  -1.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.29'
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("E38").Value = '  +1.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.99'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.09'
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '432.43'
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.73'
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.860.83'
$ws.Range("E44").Value = '  -2.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.269'
$ws.Range("E45").Value = '  -2.89%  '
$ws.Range("E46").Value = '  -2.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.18'
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.85'
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.109'
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.03'
$ws.Range("E51").Value = '  -2.90%  '
